$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.404345
$ws.Range("H2").Value = 40.213035
$ws.Range("I2").Value = 0.01122005832922476
$ws.Range("J2").Value = 0.01122005832922476
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1465046666666667
$ws.Range("N2").Value = 0.439514
$ws.Range("O2").Value = 0.07745172725947863
$ws.Range("P2").Value = 0.07745172725947864
$ws.Range("Q2").Value = 1.96379909611
$ws.Range("R2").Value = 17.67419186499
$ws.Range("S2").Value = 0.0008690128975505575
$ws.Range("T2").Value = 0.0008690128975505576
$ws.Range("G3").Value = 13.404345
$ws.Range("H3").Value = 40.213035
$ws.Range("I3").Value = 0.01122005832922476
$ws.Range("J3").Value = 0.01122005832922476
$ws.Range("N3").Value = 4.707498
$ws.Range("O3").Value = 0.8295614045753745
$ws.Range("P3").Value = 0.8295614045753745
$ws.Range("Q3").Value = 21.03364242627
$ws.Range("R3").Value = 189.30278183643
$ws.Range("S3").Value = 0.009307727347009321
$ws.Range("T3").Value = 0.009307727347009321
$ws.Range("G4").Value = 13.404345
$ws.Range("H4").Value = 40.213035
$ws.Range("I4").Value = 0.01122005832922476
$ws.Range("J4").Value = 0.01122005832922476
$ws.Range("O4").Value = 0.09298686816514685
$ws.Range("P4").Value = 0.09298686816514684
$ws.Range("Q4").Value = 2.357694710165
$ws.Range("R4").Value = 21.219252391485
$ws.Range("S4").Value = 0.00104331808466488
$ws.Range("T4").Value = 0.00104331808466488
$ws.Range("I5").Value = 0.9315566574535661
$ws.Range("J5").Value = 0.9315566574535661
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1465046666666667
$ws.Range("N5").Value = 0.439514
$ws.Range("O5").Value = 0.07745172725947863
$ws.Range("P5").Value = 0.07745172725947864
$ws.Range("Q5").Value = 163.0464003130514
$ws.Range("R5").Value = 1467.417602817462
$ws.Range("S5").Value = 0.07215067215984516
$ws.Range("T5").Value = 0.07215067215984518
$ws.Range("I6").Value = 0.9315566574535661
$ws.Range("J6").Value = 0.9315566574535661
$ws.Range("N6").Value = 4.707498
$ws.Range("O6").Value = 0.8295614045753745
$ws.Range("P6").Value = 0.8295614045753745
$ws.Range("S6").Value = 0.7727834491987213
$ws.Range("T6").Value = 0.7727834491987213
$ws.Range("I7").Value = 0.9315566574535661
$ws.Range("J7").Value = 0.9315566574535661
$ws.Range("O7").Value = 0.09298686816514685
$ws.Range("P7").Value = 0.09298686816514684
$ws.Range("S7").Value = 0.08662253609499962
$ws.Range("T7").Value = 0.08662253609499961
$ws.Range("I8").Value = 0.05722328421720919
$ws.Range("J8").Value = 0.05722328421720919
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1465046666666667
$ws.Range("N8").Value = 0.439514
$ws.Range("O8").Value = 0.07745172725947863
$ws.Range("P8").Value = 0.07745172725947864
$ws.Range("Q8").Value = 10.01554809474556
$ws.Range("R8").Value = 90.13993285271
$ws.Range("S8").Value = 0.004432042202082914
$ws.Range("T8").Value = 0.004432042202082915
$ws.Range("I9").Value = 0.05722328421720919
$ws.Range("J9").Value = 0.05722328421720919
$ws.Range("N9").Value = 4.707498
$ws.Range("O9").Value = 0.8295614045753745
$ws.Range("P9").Value = 0.8295614045753745
$ws.Range("R9").Value = 965.46083543247
$ws.Range("S9").Value = 0.04747022802964392
$ws.Range("T9").Value = 0.04747022802964392
$ws.Range("I10").Value = 0.05722328421720919
$ws.Range("J10").Value = 0.05722328421720919
$ws.Range("O10").Value = 0.09298686816514685
$ws.Range("P10").Value = 0.09298686816514684
$ws.Range("S10").Value = 0.005321013985482359
$ws.Range("T10").Value = 0.005321013985482358
